$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells F1:H1 (outlier-detection columns)
$ws.Range("F1").Value = "KNN_Outliers_MAD"
$ws.Range("G1").Value = "SVM_Outliers_MAD"
$ws.Range("H1").Value = "RF_Outliers_MAD"

# Match the header style already used by A1:E1
$ws.Range("A1").Copy()
$ws.Range("F1:H1").PasteSpecial(-4122)

# Populate boolean outlier-flag columns F:H for data rows 2-21
$ws.Cells.Item(2, 6).Value = $false
$ws.Cells.Item(2, 7).Value = $false
$ws.Cells.Item(2, 8).Value = $false
$ws.Cells.Item(3, 6).Value = $false
$ws.Cells.Item(3, 7).Value = $false
$ws.Cells.Item(3, 8).Value = $false
$ws.Cells.Item(4, 6).Value = $false
$ws.Cells.Item(4, 7).Value = $false
$ws.Cells.Item(4, 8).Value = $false
$ws.Cells.Item(5, 6).Value = $false
$ws.Cells.Item(5, 7).Value = $false
$ws.Cells.Item(5, 8).Value = $false
$ws.Cells.Item(6, 6).Value = $false
$ws.Cells.Item(6, 7).Value = $false
$ws.Cells.Item(6, 8).Value = $false
$ws.Cells.Item(7, 6).Value = $false
$ws.Cells.Item(7, 7).Value = $false
$ws.Cells.Item(7, 8).Value = $false
$ws.Cells.Item(8, 6).Value = $false
$ws.Cells.Item(8, 7).Value = $false
$ws.Cells.Item(8, 8).Value = $false
$ws.Cells.Item(9, 6).Value = $false
$ws.Cells.Item(9, 7).Value = $false
$ws.Cells.Item(9, 8).Value = $false
$ws.Cells.Item(10, 6).Value = $true
$ws.Cells.Item(10, 7).Value = $false
$ws.Cells.Item(10, 8).Value = $false
$ws.Cells.Item(11, 6).Value = $true
$ws.Cells.Item(11, 7).Value = $false
$ws.Cells.Item(11, 8).Value = $false
$ws.Cells.Item(12, 6).Value = $false
$ws.Cells.Item(12, 7).Value = $false
$ws.Cells.Item(12, 8).Value = $false
$ws.Cells.Item(13, 6).Value = $false
$ws.Cells.Item(13, 7).Value = $false
$ws.Cells.Item(13, 8).Value = $false
$ws.Cells.Item(14, 6).Value = $false
$ws.Cells.Item(14, 7).Value = $false
$ws.Cells.Item(14, 8).Value = $false
$ws.Cells.Item(15, 6).Value = $false
$ws.Cells.Item(15, 7).Value = $false
$ws.Cells.Item(15, 8).Value = $false
$ws.Cells.Item(16, 6).Value = $false
$ws.Cells.Item(16, 7).Value = $false
$ws.Cells.Item(16, 8).Value = $false
$ws.Cells.Item(17, 6).Value = $false
$ws.Cells.Item(17, 7).Value = $true
$ws.Cells.Item(17, 8).Value = $true
$ws.Cells.Item(18, 6).Value = $false
$ws.Cells.Item(18, 7).Value = $false
$ws.Cells.Item(18, 8).Value = $false
$ws.Cells.Item(19, 6).Value = $false
$ws.Cells.Item(19, 7).Value = $false
$ws.Cells.Item(19, 8).Value = $false
$ws.Cells.Item(20, 6).Value = $false
$ws.Cells.Item(20, 7).Value = $false
$ws.Cells.Item(20, 8).Value = $false
$ws.Cells.Item(21, 6).Value = $false
$ws.Cells.Item(21, 7).Value = $false
$ws.Cells.Item(21, 8).Value = $false
